$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "5`nCS-217-05812-TBD-IDE-104A`nCS-328-05807-TBD-IDE-118A`nCS-328-05807-TBD-IDE-118B"
$ws.Range("E4").Value = "5`nCS-217-05812-TBD-IDE-104A`nCS-328-05807-TBD-IDE-118A`nCS-328-05807-TBD-IDE-118B"
$ws.Range("B5").Value = "7`nCS-219-07810-TBD-IDE-204A`nCS-407-07811-TBD-IDE-317"
$ws.Range("E5").Value = "7`nCS-219-07810-TBD-IDE-204A`nCS-407-07811-TBD-IDE-317"
$ws.Range("B7").Value = "11`nCS-231-11811-Datta-IDE-317`nCS-303-11812-Anderson-IDE-113A"
$ws.Range("E7").Value = "11`nCS-231-11811-Datta-IDE-317`nCS-303-11812-Anderson-IDE-113A"
$ws.Range("C3").Value = "4`nCS-312-04815-Bancroft-IDE-317`nCS-219-04814-TBD-IDE-118A`nCS-219-04814-TBD-IDE-118B`nCS-217-04809-TBD-TBD-TBD"
$ws.Range("C4").Value = "13`nCS-113-13850-Kim-IDE-104A`nCS-113-13851-TBD-TBD-TBD`nCS-113-13852-TBD-IDE-204A"
$ws.Range("C6").Value = "10`nCS-203-10814-Datta-IDE-217A`nCS-361-10813-TBD-IDE-318"
$ws.Range("D4").Value = "6`nCS-113-06853-Kim-IDE-113A`nCS-110-06814-Bancroft-IDE-317`nCS-114-06856-TBD-IDE-104A`nCS-217-06800-TBD-IDE-217A`nCS-113-06855-TBD-TBD-TBD"
$ws.Range("D5").Value = "8`nCS-113-06853-Kim-IDE-113A`nCS-110-08816-Bancroft-IDE-317`nCS-114-06856-TBD-IDE-104A`nCS-113-06855-TBD-TBD-TBD`nCS-231-08807-TBD-TBD-TBD"
$ws.Range("D7").Value = "12`nCS-110-12809-TBD-IDE-323"
$ws.Range("F3").Value = "4.0`nCS-312-04815-Bancroft-IDE-317`nCS-219-04814-TBD-IDE-118A`nCS-219-04814-TBD-IDE-118B`nCS-217-04809-TBD-TBD-TBD"
$ws.Range("F4").Value = "6.0`nCS-113L-06854-Kim-IDE-113A`nCS-110-06814-Bancroft-IDE-317`nCS-114L-06857-TBD-IDE-104A`nCS-217-06800-TBD-IDE-217A"
$ws.Range("F5").Value = "8.0`nCS-113L-08850-Kim-IDE-104A`nCS-110-08816-Bancroft-IDE-317`nCS-113L-08851-TBD-TBD-TBD`nCS-231-08807-TBD-TBD-TBD"
$ws.Range("F6").Value = "10.0`nCS-203-10814-Datta-IDE-217A`nCS-361-10813-TBD-IDE-318`nCS-113L-10852-TBD-IDE-204A`nCS-113L-10855-TBD-TBD-TBD"
$ws.Range("F7").Value = "12.0`nCS-110-12809-TBD-IDE-323"
